$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This workbook stores one match-result per row (columns B..AD), with column A
# holding a running "id" and column D holding the match Date. The update swaps
# the match data between two pairs of rows while leaving each row's id/date in
# place (the two matches simply trade places in the result set):
#   - rows 3 and 4   (id 1 and id 2)
#   - rows 117 and 118 (id 115 and id 116)
# Columns swapped: B,C,E,F,G,H,I,J,K,L,M,N,O,P,Q,R,S,T,U,V,W,X,Y,Z,AA,AB,AC,AD
# (everything except A = id and D = Date).

function Swap-MatchRows {
    param($Row1, $Row2)

    $range1 = $ws.Range("B$Row1`:C$Row1")
    $range2 = $ws.Range("B$Row2`:C$Row2")
    $vals1 = $range1.Value2
    $vals2 = $range2.Value2
    $range1.Value2 = $vals2
    $range2.Value2 = $vals1

    $range1 = $ws.Range("E$Row1`:AD$Row1")
    $range2 = $ws.Range("E$Row2`:AD$Row2")
    $vals1 = $range1.Value2
    $vals2 = $range2.Value2
    $range1.Value2 = $vals2
    $range2.Value2 = $vals1
}

Swap-MatchRows 3 4
Swap-MatchRows 117 118
